# Second version of IAV model without variable D and with some modifications.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# Sheet1: update the raw "Lung CD8" / "Lung CD8_E" source percentages (day 8 /
# day 10, ZT23 / ZT11) that feed the existing B79:C80 / B84:C85 formulas.
# ---------------------------------------------------------------------------
$ws1.Range("B70").Value2 = 22.666667
$ws1.Range("C70").Value2 = 19.28
$ws1.Range("B71").Value2 = 30
$ws1.Range("C71").Value2 = 28.774999999999999

$ws1.Range("B74").Value2 = 18.666666666666668
$ws1.Range("C74").Value2 = 15.68
$ws1.Range("B75").Value2 = 25.241667
$ws1.Range("C75").Value2 = 23.4

# The existing formulas in B79:C80 and B84:C85 recompute automatically from
# the new inputs above, so nothing further is required there.

# ---------------------------------------------------------------------------
# New block: "Lung CD8 Tcell (as % of all live cells) without T_E"
# (T cell minus T_E cell, per day) in rows 87-90.
# ---------------------------------------------------------------------------
$ws1.Range("A87").Value = "Lung CD8 Tcell (as % of all live cells) without T_E"

$ws1.Range("B88").Value = "ZT23"
$ws1.Range("C88").Value = "ZT11"

$ws1.Range("A89").Value = "day8"
$ws1.Range("B89").Formula = "=B79-B84"
$ws1.Range("C89").Formula = "=C79-C84"

$ws1.Range("A90").Value = "day10"
$ws1.Range("B90").Formula = "=B80-B85"
$ws1.Range("C90").Formula = "=C80-C85"

# ---------------------------------------------------------------------------
# Sheet2: second table (rows 13-22) gains the same 24/48/96/144 hr columns as
# the first table (rows 1-10), plus the day8/day10 T & T_E values shrink to
# the new per-timepoint series.
# ---------------------------------------------------------------------------
$ws2.Range("B20").Value2 = 24
$ws2.Range("C20").Value2 = 48
$ws2.Range("D20").Value2 = 96
$ws2.Range("E20").Value2 = 144
$ws2.Range("F20").Formula = "=8*24"
$ws2.Range("G20").Formula = "=10*24"

$ws2.Range("B21").Value2 = 0.05
$ws2.Range("C21").Value2 = 0.08
$ws2.Range("D21").Value2 = 0.1
$ws2.Range("E21").Value2 = 0.5
$ws2.Range("F21").Value2 = 2.3328000000000002
$ws2.Range("G21").Value2 = 3.4830000000000001

$ws2.Range("B22").Value2 = 0.05
$ws2.Range("C22").Value2 = 0.1
$ws2.Range("D22").Value2 = 0.5
$ws2.Range("E22").Value2 = 1
$ws2.Range("F22").Value2 = 10.160639999999999
$ws2.Range("G22").Value2 = 15.1632

# ---------------------------------------------------------------------------
# Selections: Sheet1 ends up focused on C85, Sheet2 (the active tab) on B22.
# ---------------------------------------------------------------------------
$ws1.Range("C85").Select()
$ws2.Activate()
$ws2.Range("B22").Select()
